$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-21 Friday" "2024-06-22 Saturday"
Replace-Text "562×6=" "870×5="
Replace-Text "139×6=" "278×6="
Replace-Text "872×9=" "663×6="
Replace-Text "423×6=" "614×9="
Replace-Text "739×5=" "397×9="
Replace-Text "880×6=" "933×9="
Replace-Text "183×4=" "292×4="
Replace-Text "294×9=" "780×6="
Replace-Text "755×6=" "838×8="
Replace-Text "837×2=" "644×6="
Replace-Text "663×4=" "230×5="
Replace-Text "634×8=" "431×9="
Replace-Text "611×3=" "415×7="
Replace-Text "828×4=" "654×3="
Replace-Text "958×2=" "359×5="
Replace-Text "296×4=" "521×6="
Replace-Text "165×8=" "963×7="
Replace-Text "719×5=" "317×7="
Replace-Text "551×9=" "504×2="
Replace-Text "664×7=" "991×2="
Replace-Text "347×8=" "714×9="
Replace-Text "649×7=" "819×8="
Replace-Text "569×5=" "581×6="
Replace-Text "240×3=" "522×9="
Replace-Text "313×2=" "609×6="
